$wb = $excel.ActiveWorkbook

# --- Rename existing headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row (match the bold/bordered/centered look used by the other sheets' headers) ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Forecast data rows ---
$data = @(
    @(45025.99999999999, 4, -1.465047091100585, 8.721349448978684),
    @(45039.99999999999, 4, -1.218164305656335, 8.625001418188198),
    @(45060.99999999999, 4, -1.01804076680148, 8.927092067107399),
    @(45067.99999999999, 4, -0.9900171378953828, 9.172845771710323),
    @(45074.99999999999, 4, -1.351191174534009, 9.194886776409616),
    @(45081.99999999999, 4, -1.04060641512677, 9.318777761944245),
    @(45088.99999999999, 4, -0.78077949216098, 9.251523822504959),
    @(45095.99999999999, 4, -1.168876378835268, 9.253647348706609),
    @(45102.99999999999, 4, -1.081521224632676, 8.987879919787845),
    @(45109.99999999999, 4, -1.025940352520969, 9.244892600507168),
    @(45116.99999999999, 4, -1.062999375643345, 9.612285268304145),
    @(45123.99999999999, 4, -0.55495387297064, 9.268608687569362),
    @(45130.99999999999, 4, -0.7702631572198879, 9.248348888377626),
    @(45137.99999999999, 4, -0.8080924491912328, 9.612903251837633),
    @(45144.99999999999, 4, -0.8016732486640376, 9.629051059997119),
    @(45158.99999999999, 4, -0.4958653362484319, 9.705307688917165),
    @(45172.99999999999, 5, -0.8056943643503841, 9.986534939633557),
    @(45179.99999999999, 5, -0.8737817424715502, 9.524553193613354),
    @(45186.99999999999, 5, -0.4913655528375346, 9.822134737033329),
    @(45200.99999999999, 5, -0.3319834905668708, 9.974806314997309),
    @(45214.99999999999, 5, -0.5217275112913795, 9.66501753804358),
    @(45235.99999999999, 5, -0.2174703840746353, 10.08708428910164),
    @(45242.99999999999, 5, -0.4006302943044324, 10.13517641306772),
    @(45249.99999999999, 5, -0.5308600339106324, 10.23443885149865),
    @(45312.99999999999, 5, 0.5781504481159025, 10.36524025427234),
    @(45333.99999999999, 5, 0.5937852114453425, 10.26390399272246),
    @(45347.99999999999, 5, 0.4523599356231421, 10.69906202253646),
    @(45354.99999999999, 6, 0.3030563720208921, 10.62497177155022),
    @(45361.99999999999, 6, 0.4372757123584119, 10.37796343552764),
    @(45368.99999999999, 6, 0.4228230483939389, 10.74650401405752),
    @(45375.99999999999, 6, -0.07740841471863029, 10.99314531080965),
    @(45382.99999999999, 6, 0.611910875811399, 10.84254502623224),
    @(45445.99999999999, 6, 0.6082174341513537, 11.33324047520678),
    @(45452.99999999999, 6, 1.041656044709551, 10.96369317452782),
    @(45459.99999999999, 6, 0.6328748998815669, 11.17911563923823),
    @(45466.99999999999, 6, 0.8520808473485088, 11.20073191266569),
    @(45480.99999999999, 6, 0.9361946908026789, 11.24634612896548),
    @(45487.99999999999, 6, 1.051645980356215, 11.514520408398),
    @(45494.99999999999, 6, 0.7802262421137881, 11.30865523763319),
    @(45536.99999999999, 7, 1.316691131025918, 11.87239627363113),
    @(45550.99999999999, 7, 1.426147775412971, 11.65585248234031),
    @(45557.99999999999, 7, 1.734500294860874, 11.72536675531176),
    @(45564.99999999999, 7, 1.302928168755259, 12.08571809683603),
    @(45571.99999999999, 7, 1.76683662131716, 12.22635275699876),
    @(45578.99999999999, 7, 1.903092806050049, 11.61692776415298),
    @(45585.99999999999, 7, 1.62548034270361, 12.05904486838086),
    @(45592.99999999999, 7, 1.768310401034605, 12.21229600812588),
    @(45599.99999999999, 7, 2.022860878857797, 12.19034707586148),
    @(45606.99999999999, 7, 1.954442681885999, 11.66021418720003),
    @(45613.99999999999, 7, 2.228444714654342, 12.4392385851361),
    @(45620.99999999999, 7, 1.573883561082443, 11.78525000976312)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the date/time display format used by the other "Order Week"/"Order Month" columns
$wsForecast.Range("A2:A52").NumberFormat = $wsWeekly.Range("A2").NumberFormat
